$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.300.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.648.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.95%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("E9").Value = "  +2.74%  "

$ws.Range("E10").Value = "  -1.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "

$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.131.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.281.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.643.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "364.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.85%  "

$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.66%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  -2.81%  "

$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "558.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.72%  "

$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("E32").Value = "  -1.61%  "

$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +0.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("E39").Value = "  +1.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("E41").Value = "  -0.76%  "

$ws.Range("E42").Value = "  +4.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "159.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.14%  "

$ws.Range("E48").Value = "  -0.99%  "

$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  -0.35%  "

$ws.Range("E51").Value = "  +0.85%  "
